# Update the "Förändrad" (Changed) date column (C2:C84) from 46081 to 46082
# (i.e. advance the date serial number by one day) for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value = 46082
    }
}
